# Rename the existing sheet, then build two additional sheets ("中坚组" and
# "挑战组") that contain the same header/data rows (A1:Z5) as the renamed
# sheet, each with their own view/column-width settings, matching the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "传奇组"

# --- create "中坚组" right after the first sheet ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "中坚组"

$ws1.Range("A1:Z5").Copy()
$ws2.Range("A1").PasteSpecial(-4104)

$ws2.Columns.Item(1).ColumnWidth = 13.7109375
$ws2.Range("A1:XFD5").Select()

# --- create "挑战组" right after "中坚组" --------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "挑战组"

$ws1.Range("A1:Z5").Copy()
$ws3.Range("A1").PasteSpecial(-4104)

$ws3.Columns.Item(1).ColumnWidth = 19.85546875
$ws3.Range("B4").Select()

# --- tidy up the first sheet's selection/view ---------------------------
$ws1.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 1

$ws1.Activate()
